$d = $word.ActiveDocument

function Get-ParaIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -eq ($text + "`r")) {
            return $i
        }
    }
    return -1
}

function Set-ParaTextSplit($idx, $parts) {
    # Replace paragraph $idx's content with multiple runs built from $parts
    # (array of strings), by typing the first chunk and then appending the
    # rest one at a time.
    $p = $d.Paragraphs.Item($idx)
    $r = $p.Range
    $r.Text = $parts[0]
    for ($j = 1; $j -lt $parts.Count; $j++) {
        $ip = $d.Range($r.Start, $r.End - 1)
        $ip.InsertAfter($parts[$j])
    }
}

# ---------------------------------------------------------------------------
# Step 1: "Were there any points where you felt unable to do what you
# wanted?" -> "What minor bugs and glitches did you encounter?"
# ---------------------------------------------------------------------------
$idx = Get-ParaIndex("Were there any points where you felt unable to do what you wanted?")
Set-ParaTextSplit $idx @("What", " minor", " bugs and glitches did you encounter", "?")

# ---------------------------------------------------------------------------
# Step 2: insert a new paragraph right after it:
# "What game-breaking bugs and glitches did you encounter?"
# ---------------------------------------------------------------------------
$idx = Get-ParaIndex("What minor bugs and glitches did you encounter?")
$d.Paragraphs.Item($idx).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($idx + 1).Range.Text = "What game-breaking bugs and glitches did you encounter?"

# ---------------------------------------------------------------------------
# Step 3: drop the now-redundant "What bugs and glitches did you encounter?"
# and "How would you describe the game to friends/family?" paragraphs.
# ---------------------------------------------------------------------------
$idx = Get-ParaIndex("What bugs and glitches did you encounter?")
$d.Paragraphs.Item($idx).Range.Delete() | Out-Null

$idx = Get-ParaIndex("How would you describe the game to friends/family?")
$d.Paragraphs.Item($idx).Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# Step 4: "Which weapons did you prioritise using? Why?" ->
# "Which weapons did you prioritise using?" (the "Why?" resurfaces below).
# ---------------------------------------------------------------------------
$idx = Get-ParaIndex("Which weapons did you prioritise using? Why?")
Set-ParaTextSplit $idx @("Which weapons did you prioritise using", "?")

# ---------------------------------------------------------------------------
# Step 5: relocate "What would you add, remove, or change about the
# experience?" to just before the final question, appending " Why?".
# ---------------------------------------------------------------------------
$idx = Get-ParaIndex("What would you add, remove, or change about the experience?")
$d.Paragraphs.Item($idx).Range.Delete() | Out-Null

$idx = Get-ParaIndex("Is there anything else you would like to say about the game?")
$d.Paragraphs.Item($idx).Range.InsertParagraphBefore() | Out-Null
Set-ParaTextSplit $idx @("What would you add, remove, or change about the experience?", " Why?")

Write-Output "Paragraphs count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "$i : [$($p.Range.Text)]"
}
